# Add new lookup groups (INVENTORY_TYPE, STOCK_TYPE, INVENTORY_STATUS,
# SUBINVENTORIES, COSTING_METHOD) to the Lookup_Values sheet, appended
# after the existing data (which ends at row 116).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

function Add-LookupRow {
    param(
        [int]$RowNum,
        [string]$LookupType,
        [string]$Value,
        [string]$Description,
        [bool]$KeepStyle
    )

    # Inserting a whole row (shift down) makes this engine copy the
    # formatting of the row immediately above onto the new row - which
    # gives column B the same style used throughout the sheet.
    $ws.Rows.Item($RowNum).Insert(-4121)

    $ws.Range("A$RowNum").Value = $LookupType
    $ws.Range("B$RowNum").Value = $Value
    $ws.Range("C$RowNum").Value = $Description

    if (-not $KeepStyle) {
        # Several of the appended rows were typed in without inheriting
        # the usual formatting - strip it back to the default style.
        $ws.Range("A$RowNum`:C$RowNum").ClearFormats()
    }
}

Add-LookupRow 117 "INVENTORY_TYPE"   "Inventory Items"      "That are stocked and tracked in inventory e.g raw materials, finished goods, or work-in-progress." $true
Add-LookupRow 118 "INVENTORY_TYPE"   "Non-Inventory Items"  "That are not stocked or tracked in inventory. They are typically used for services or items that are expensed immediately upon purchase." $true
Add-LookupRow 119 "INVENTORY_TYPE"   "Kit Items"            "These are items that are made up of a collection of other items." $true

Add-LookupRow 120 "STOCK_TYPE"       "On-Hand Quantity"     "This refers to the actual quantity of items physically present in the inventory." $true
Add-LookupRow 121 "STOCK_TYPE"       "Reserved Quantity"    "This is the quantity of items that have been reserved for specific sales orders, work orders, or other transactions." $true
Add-LookupRow 122 "STOCK_TYPE"       "Available Quantity"   "This is the quantity of items that are available for use, calculated as On-Hand Quantity minus Reserved Quantity." $true

Add-LookupRow 123 "INVENTORY_STATUS" "Active"   "Items that are currently available for use in transactions." $false
Add-LookupRow 124 "INVENTORY_STATUS" "Inactive" "Items that are no longer available for use but may still exist in the system for historical purposes." $false
Add-LookupRow 125 "INVENTORY_STATUS" "Obsolete" "Items that are no longer used and are typically removed from active inventory." $false

Add-LookupRow 126 "SUBINVENTORIES" "Stores"                  "Subinventories used for storing raw materials and components." $false
Add-LookupRow 127 "SUBINVENTORIES" "Finished Goods"          "Subinventories used for storing completed products ready for sale." $false
Add-LookupRow 128 "SUBINVENTORIES" "WIP (Work in Progress)"  "Subinventories used for items that are in the process of being manufactured." $false

Add-LookupRow 129 "COSTING_METHOD" "Standard Costing"            "Inventory is valued at a predetermined standard cost." $false
Add-LookupRow 130 "COSTING_METHOD" "Average Costing"             "Inventory is valued at an average cost, which is recalculated after each transaction." $false
Add-LookupRow 131 "COSTING_METHOD" "FIFO (First In, First Out)"  "Inventory is valued based on the assumption that the oldest items are sold first." $false
Add-LookupRow 132 "COSTING_METHOD" "LIFO (Last In, First Out)"   "Inventory is valued based on the assumption that the newest items are sold first." $false

# Column C now holds much longer descriptions - widen it to fit like Excel's
# "AutoFit Column Width" would.
$ws.Columns.Item(3).AutoFit()

# Reflect where the user ended up after typing all this in: scrolled down
# with the next empty row selected.
$ws.Range("A136").Select()

Write-Host "Added 16 new lookup rows (117-132)."
